# Generate Report for Handoff
# Adds the new file "c212e0aa-1431-4244-a13f-ef940ede77e2" as row 4 on all
# three worksheets (Overview, zh-cn, de-de), matching the existing pattern
# used for the other two files already on the sheets.

$wb = $excel.ActiveWorkbook

$newFile        = "c212e0aa-1431-4244-a13f-ef940ede77e2"
$mdName         = "$newFile.md"
$zhTargetFile   = "$newFile.48923b1f5155f022de314ae74376fccf82795e91.zh-cn.xlf"
$deTargetFile   = "$newFile.48923b1f5155f022de314ae74376fccf82795e91.de-de.xlf"
$mdUrl          = "https://github.com/OpenLocalizationTest/oltest/blob/e0738bf1070e4f694726817815c2060ae43284eb/e2e/$mdName"
$zhTargetUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aabd247eff1ca40c968caba4c1a81eac4a0426f5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhTargetFile"
$deTargetUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/113c2e237e9628c492ebf5eb614b3c8875844b05/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deTargetFile"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$rowOv = 4
$wsOverview.Cells.Item($rowOv, 1).Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rowOv, 1), $mdUrl, "", "", $mdName)
$wsOverview.Cells.Item($rowOv, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item($rowOv, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item($rowOv, 4).Value = "2016-55-13 20:55:43"

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File
# | Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$rowZh = 4
$wsZh.Cells.Item($rowZh, 1).Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rowZh, 1), $mdUrl, "", "", $mdName)
$wsZh.Cells.Item($rowZh, 2).Value = ".md"
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rowZh, 2), $mdUrl, "", "", ".md")
$wsZh.Cells.Item($rowZh, 3).Value = "Ready for handoff"
$wsZh.Cells.Item($rowZh, 4).Value = $zhTargetFile
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rowZh, 4), $zhTargetUrl, "", "", $zhTargetFile)
$wsZh.Cells.Item($rowZh, 5).Value = "2016-03-13 20:55:39"
$wsZh.Cells.Item($rowZh, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item($rowZh, 8).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($rowZh, 9).Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet: same shape as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$rowDe = 4
$wsDe.Cells.Item($rowDe, 1).Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rowDe, 1), $mdUrl, "", "", $mdName)
$wsDe.Cells.Item($rowDe, 2).Value = ".md"
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rowDe, 2), $mdUrl, "", "", ".md")
$wsDe.Cells.Item($rowDe, 3).Value = "Ready for handoff"
$wsDe.Cells.Item($rowDe, 4).Value = $deTargetFile
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rowDe, 4), $deTargetUrl, "", "", $deTargetFile)
$wsDe.Cells.Item($rowDe, 5).Value = "2016-03-13 20:55:43"
$wsDe.Cells.Item($rowDe, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item($rowDe, 8).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($rowDe, 9).Value = "Include"
